$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5-36 of columns A:F (doctyp_code, doccat_code, lang_code, is_active, cr_by, cr_dtimes).
# Row 5-7 are corrected/re-ordered existing rows; rows 8-36 are newly added
# Document Type records (plus the extra Mac-Address / POx category rows).
$data = @(
  @("DOC001", "POI", "ara", $true,  "superadmin", "now()"),
  @("CRN",    "POR", "ara", $true,  "superadmin", "now()"),
  @("COB",    "POB", "ara", $false, "superadmin", "now()"),
  @("DOC001", "POI", "ara", $true,  "superadmin", "now()"),
  @("DOC002", "POI", "ara", $true,  "superadmin", "now()"),
  @("DOC003", "POI", "ara", $true,  "superadmin", "now()"),
  @("DOC004", "POI", "ara", $true,  "superadmin", "now()"),
  @("DOC005", "POI", "ara", $true,  "superadmin", "now()"),
  @("DOC006", "POI", "ara", $true,  "superadmin", "now()"),
  @("DOC007", "POI", "ara", $true,  "superadmin", "now()"),
  @("DOC008", "POI", "ara", $true,  "superadmin", "now()"),
  @("DOC009", "POI", "ara", $true,  "superadmin", "now()"),
  @("DOC010", "POI", "ara", $true,  "superadmin", "now()"),
  @("DOC011", "POI", "ara", $true,  "superadmin", "now()"),
  @("DOC012", "POI", "ara", $true,  "superadmin", "now()"),
  @("DOC001", "POA", "ara", $true,  "superadmin", "now()"),
  @("DOC013", "POA", "ara", $true,  "superadmin", "now()"),
  @("DOC014", "POA", "ara", $true,  "superadmin", "now()"),
  @("DOC015", "POA", "ara", $true,  "superadmin", "now()"),
  @("DOC004", "POA", "ara", $true,  "superadmin", "now()"),
  @("DOC005", "POA", "ara", $true,  "superadmin", "now()"),
  @("DOC006", "POA", "ara", $true,  "superadmin", "now()"),
  @("DOC016", "POA", "ara", $true,  "superadmin", "now()"),
  @("DOC017", "POA", "ara", $true,  "superadmin", "now()"),
  @("DOC018", "POA", "ara", $true,  "superadmin", "now()"),
  @("DOC008", "POA", "ara", $true,  "superadmin", "now()"),
  @("DOC024", "POR", "ara", $true,  "superadmin", "now()"),
  @("DOC025", "POR", "ara", $true,  "superadmin", "now()"),
  @("DOC026", "POR", "ara", $true,  "superadmin", "now()"),
  @("DOC001", "POR", "ara", $true,  "superadmin", "now()"),
  @("DOC027", "POR", "ara", $true,  "superadmin", "now()"),
  @("DOC028", "POR", "ara", $true,  "superadmin", "now()")
)

$rowCount = $data.Count
$colCount = 6

$arr = New-Object 'object[,]' $rowCount, $colCount
for ($r = 0; $r -lt $rowCount; $r++) {
  for ($c = 0; $c -lt $colCount; $c++) {
    $arr[$r, $c] = $data[$r][$c]
  }
}

$startRow = 5
$endRow = $startRow + $rowCount - 1
$ws.Range("A" + $startRow + ":F" + $endRow).Value = $arr

# Keep the header-row "select everything below the table" selection in sync
# with the new table size (previously A8:XFD1048576, now starts one column
# to the right of the last used column and one row below the header).
$ws.Range("G1:XFD1048576").Select()
